# Fix sample company data: proper expense products and categories.
#
# 1) Expenses sheet: rename the "(N units)" product-description suffixes to
#    "(Inventory)" so they differentiate from the revenue products that share
#    the same base name.
# 2) Products sheet: append 20 new expense-type products (PRD-021..PRD-040)
#    covering Inventory Purchases, Office Supplies, Software Subscriptions,
#    Professional Services and Marketing categories.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Expenses sheet - rename "(N units)" suffix to "(Inventory)"
# ---------------------------------------------------------------------
$wsExpenses = $wb.Worksheets.Item("Expenses")

$wsExpenses.Range("D2").Value  = "Dell Latitude 5540 Laptop (Inventory)"
$wsExpenses.Range("D3").Value  = "Dell OptiPlex 7010 Desktop (Inventory)"
$wsExpenses.Range("D4").Value  = "Dell UltraSharp U2722D Monitor (Inventory)"
$wsExpenses.Range("D5").Value  = "Lenovo ThinkPad T14 (Inventory)"
$wsExpenses.Range("D6").Value  = "Cisco Meraki MX68 Firewall (Inventory)"
$wsExpenses.Range("D7").Value  = "Dell P2422H Monitor (Inventory)"
$wsExpenses.Range("D8").Value  = "Dell Precision 3660 Workstation (Inventory)"
$wsExpenses.Range("D9").Value  = "Cisco Catalyst 1000-24T Switch (Inventory)"
$wsExpenses.Range("D10").Value = "Dell Latitude 7440 Laptop (Inventory)"
$wsExpenses.Range("D12").Value = "Dell Latitude 5540 Laptop (Inventory)"
$wsExpenses.Range("D13").Value = "Dell OptiPlex 7010 Desktop (Inventory)"
$wsExpenses.Range("D14").Value = "Ubiquiti UniFi Access Point (Inventory)"
$wsExpenses.Range("D15").Value = "Dell Latitude 5540 Laptop (Inventory)"
$wsExpenses.Range("D16").Value = "Dell Latitude 7440 Laptop (Inventory)"
$wsExpenses.Range("D17").Value = "Dell Precision 3660 Workstation (Inventory)"

# ---------------------------------------------------------------------
# 2) Products sheet - add PRD-021..PRD-040 (expense-type products)
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")

$newProducts = @(
    @("PRD-021","Dell Latitude 5540 Laptop (Inventory)","Expenses","Product","EXP-DELL-LAT5540","Purchase of Dell Latitude 5540 laptops for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",5,50),
    @("PRD-022","Dell Latitude 7440 Laptop (Inventory)","Expenses","Product","EXP-DELL-LAT7440","Purchase of Dell Latitude 7440 laptops for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",5,50),
    @("PRD-023","Dell OptiPlex 7010 Desktop (Inventory)","Expenses","Product","EXP-DELL-OPT7010","Purchase of Dell OptiPlex desktops for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",5,50),
    @("PRD-024","Dell Precision 3660 Workstation (Inventory)","Expenses","Product","EXP-DELL-PREC3660","Purchase of Dell Precision workstations for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",3,25),
    @("PRD-025","Dell UltraSharp U2722D Monitor (Inventory)","Expenses","Product","EXP-DELL-U2722D","Purchase of Dell UltraSharp monitors for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",5,50),
    @("PRD-026","Dell P2422H Monitor (Inventory)","Expenses","Product","EXP-DELL-P2422H","Purchase of Dell P2422H monitors for inventory","CAT-PUR-001","Inventory Purchases","SUP-001","Dell Technologies",5,50),
    @("PRD-027","Lenovo ThinkPad T14 (Inventory)","Expenses","Product","EXP-LEN-T14","Purchase of Lenovo ThinkPad laptops for inventory","CAT-PUR-001","Inventory Purchases","SUP-007","Lenovo",5,50),
    @("PRD-028","Cisco Meraki MX68 Firewall (Inventory)","Expenses","Product","EXP-CISCO-MX68","Purchase of Cisco Meraki firewalls for inventory","CAT-PUR-001","Inventory Purchases","SUP-004","Cisco Systems",2,20),
    @("PRD-029","Cisco Catalyst 1000-24T Switch (Inventory)","Expenses","Product","EXP-CISCO-CAT1000","Purchase of Cisco switches for inventory","CAT-PUR-001","Inventory Purchases","SUP-004","Cisco Systems",2,20),
    @("PRD-030","Ubiquiti UniFi Access Point (Inventory)","Expenses","Product","EXP-UBNT-UAP","Purchase of Ubiquiti access points for inventory","CAT-PUR-001","Inventory Purchases","SUP-002","Ingram Micro",5,50),
    @("PRD-031","Office Supplies - Q1","Expenses","Product","EXP-OFF-Q1","Q1 office supplies purchase","CAT-PUR-002","Office Supplies","SUP-002","Ingram Micro",0,0),
    @("PRD-032","Office Supplies - Q2","Expenses","Product","EXP-OFF-Q2","Q2 office supplies purchase","CAT-PUR-002","Office Supplies","SUP-002","Ingram Micro",0,0),
    @("PRD-033","Office Supplies - Q3","Expenses","Product","EXP-OFF-Q3","Q3 office supplies purchase","CAT-PUR-002","Office Supplies","SUP-002","Ingram Micro",0,0),
    @("PRD-034","Office Supplies - Q4","Expenses","Product","EXP-OFF-Q4","Q4 office supplies purchase","CAT-PUR-002","Office Supplies","SUP-002","Ingram Micro",0,0),
    @("PRD-035","Logitech Peripherals Bundle","Expenses","Product","EXP-LOG-BUNDLE","Logitech keyboard, mouse, and webcam bundle","CAT-PUR-001","Inventory Purchases","SUP-002","Ingram Micro",5,50),
    @("PRD-036","Adobe Creative Cloud (8 licenses)","Expenses","Service","EXP-ADOBE-CC8","Adobe Creative Cloud subscription (8 licenses)","CAT-PUR-003","Software Subscriptions","SUP-005","CDW Corporation",0,0),
    @("PRD-037","Microsoft 365 Business (40 licenses)","Expenses","Service","EXP-MS365-40","Microsoft 365 Business subscription (40 licenses)","CAT-PUR-003","Software Subscriptions","SUP-003","Tech Data Corporation",0,0),
    @("PRD-038","Professional Development Training","Expenses","Service","EXP-TRAIN","Employee professional development and training","CAT-PUR-004","Professional Services","","",0,0),
    @("PRD-039","Marketing - Fall Campaign","Expenses","Service","EXP-MKT-FALL","Fall marketing campaign expenses","CAT-PUR-005","Marketing","","",0,0),
    @("PRD-040","Marketing - Spring Campaign","Expenses","Service","EXP-MKT-SPRING","Spring marketing campaign expenses","CAT-PUR-005","Marketing","","",0,0)
)

$row = 22
foreach ($p in $newProducts) {
    $wsProducts.Cells.Item($row, 1).Value  = $p[0]
    $wsProducts.Cells.Item($row, 2).Value  = $p[1]
    $wsProducts.Cells.Item($row, 3).Value  = $p[2]
    $wsProducts.Cells.Item($row, 4).Value  = $p[3]
    $wsProducts.Cells.Item($row, 5).Value  = $p[4]
    $wsProducts.Cells.Item($row, 6).Value  = $p[5]
    $wsProducts.Cells.Item($row, 7).Value  = $p[6]
    $wsProducts.Cells.Item($row, 8).Value  = $p[7]
    if ($p[8] -ne "") {
        $wsProducts.Cells.Item($row, 9).Value = $p[8]
    }
    if ($p[9] -ne "") {
        $wsProducts.Cells.Item($row, 10).Value = $p[9]
    }
    $wsProducts.Cells.Item($row, 11).Value = $p[10]
    $wsProducts.Cells.Item($row, 12).Value = $p[11]
    $row = $row + 1
}
